$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Login with valid username and password"
$ws.Range("B18").Value = "FAILED"
$ws.Range("C18").Value = "chrome"
$ws.Range("D18").Value = "31_08_23125220"
$ws.Range("A19").Value = "Login with valid username and password"
$ws.Range("B19").Value = "FAILED"
$ws.Range("C19").Value = "chrome"
$ws.Range("D19").Value = "31_08_23125947"
$ws.Range("A20").Value = "Login with valid username and password"
$ws.Range("B20").Value = "FAILED"
$ws.Range("C20").Value = "chrome"
$ws.Range("D20").Value = "03_12_23194215"
$ws.Range("A21").Value = "Login with valid username and password"
$ws.Range("B21").Value = "FAILED"
$ws.Range("C21").Value = "chrome"
$ws.Range("D21").Value = "03_12_23200244"
$ws.Range("A22").Value = "Login with valid username and password"
$ws.Range("B22").Value = "FAILED"
$ws.Range("C22").Value = "chrome"
$ws.Range("D22").Value = "03_12_23201143"
$ws.Range("A23").Value = "Login with valid username and password"
$ws.Range("B23").Value = "FAILED"
$ws.Range("C23").Value = "chrome"
$ws.Range("D23").Value = "03_12_23204722"
$ws.Range("A24").Value = "Add position categories"
$ws.Range("B24").Value = "FAILED"
$ws.Range("C24").Value = "chrome"
$ws.Range("D24").Value = "03_12_23204733"
$ws.Range("A25").Value = "Edit position categories"
$ws.Range("B25").Value = "FAILED"
$ws.Range("C25").Value = "chrome"
$ws.Range("D25").Value = "03_12_23204741"
$ws.Range("A26").Value = "Delete position categories"
$ws.Range("B26").Value = "FAILED"
$ws.Range("C26").Value = "chrome"
$ws.Range("D26").Value = "03_12_23204750"
$ws.Range("A27").Value = "Add attestations"
$ws.Range("B27").Value = "FAILED"
$ws.Range("C27").Value = "chrome"
$ws.Range("D27").Value = "03_12_23204756"
$ws.Range("A28").Value = "Login with valid username and password"
$ws.Range("B28").Value = "FAILED"
$ws.Range("C28").Value = "chrome"
$ws.Range("D28").Value = "03_12_23211342"
$ws.Range("A29").Value = "Login with valid username and password"
$ws.Range("B29").Value = "FAILED"
$ws.Range("C29").Value = "chrome"
$ws.Range("D29").Value = "03_12_23213607"
$ws.Range("A30").Value = "Login with valid username and password"
$ws.Range("B30").Value = "FAILED"
$ws.Range("C30").Value = "chrome"
$ws.Range("D30").Value = "03_12_23214059"
$ws.Range("A31").Value = "Login with valid username and password"
$ws.Range("B31").Value = "FAILED"
$ws.Range("C31").Value = "chrome"
$ws.Range("D31").Value = "03_12_23214443"
$ws.Range("A32").Value = "Login with valid username and password"
$ws.Range("B32").Value = "FAILED"
$ws.Range("C32").Value = "chrome"
$ws.Range("D32").Value = "03_12_23215557"
$ws.Range("A33").Value = "Login with valid username and password"
$ws.Range("B33").Value = "FAILED"
$ws.Range("C33").Value = "chrome"
$ws.Range("D33").Value = "03_12_23215709"
